$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted above the current row 295, pushing the
# existing rows 295-368 down to 296-369 (dimension grows from T368 to T369).
$ws.Rows.Item(295).Insert()

# Fill in the newly inserted row 295 with the new record's data.
$ws.Range("A295").Value = 10
$ws.Range("B295").Value = 'Vega Modelo de Temuco'
$ws.Range("C295").Value = 'La Araucanía'
$ws.Range("D295").Value = 44722
$ws.Range("E295").Value = 9
$ws.Range("F295").Value = 'Fruta'
$ws.Range("G295").Value = 100108
$ws.Range("H295").Value = 'Tropicales y subtropicales'
$ws.Range("I295").Value = 100108002
$ws.Range("J295").Value = 'Mango'
$ws.Range("K295").Value = 'Sin especificar'
$ws.Range("L295").Value = 'Primera'
$ws.Range("M295").Value = 250
$ws.Range("N295").Value = 9000
$ws.Range("O295").Value = 9000
$ws.Range("P295").Value = 9000
$ws.Range("Q295").Value = '$/bandeja 4 kilos'
$ws.Range("R295").Value = 'Brasil'
$ws.Range("S295").Value = 2250
$ws.Range("T295").Value = 4
